$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '275.78'
$ws.Range('E2').NumberFormat = "@"
$ws.Range('E2').Value = '-1.72%'

# Row 3
$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '27.24'
$ws.Range('E3').NumberFormat = "@"
$ws.Range('E3').Value = '0.76%'

# Row 4
$ws.Range('D4').NumberFormat = "@"
$ws.Range('D4').Value = '4.757'
$ws.Range('E4').NumberFormat = "@"
$ws.Range('E4').Value = '-3.70%'

# Row 5
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '0.06321'
$ws.Range('E5').NumberFormat = "@"
$ws.Range('E5').Value = '-1.01%'

# Row 6
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '6.941'
$ws.Range('E6').NumberFormat = "@"
$ws.Range('E6').Value = '-0.25%'

# Row 7
$ws.Range('B7').Value = 'FTXToken'
$ws.Range('C7').Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '1.356'
$ws.Range('E7').NumberFormat = "@"
$ws.Range('E7').Value = '42.09%'

# Row 8
$ws.Range('B8').Value = 'MXToken'
$ws.Range('C8').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.8766'
$ws.Range('E8').NumberFormat = "@"
$ws.Range('E8').Value = '-1.00%'

# Row 9
$ws.Range('B9').Value = 'WazirX'
$ws.Range('C9').Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.1519'
$ws.Range('E9').NumberFormat = "@"
$ws.Range('E9').Value = '3.11%'

# Row 10
$ws.Range('B10').Value = 'LiechtensteinCryptoassetsExchange'
$ws.Range('C10').Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '0.05017'
$ws.Range('E10').NumberFormat = "@"
$ws.Range('E10').Value = '-3.12%'

# Row 11
$ws.Range('B11').Value = 'MandalaExchangeToken'
$ws.Range('C11').Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.07472'
$ws.Range('E11').NumberFormat = "@"
$ws.Range('E11').Value = '0.82%'

# Row 12
$ws.Range('B12').Value = 'BitrueCoin'
$ws.Range('C12').Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '0.02900'
$ws.Range('E12').NumberFormat = "@"
$ws.Range('E12').Value = '-7.03%'

# Row 13
$ws.Range('B13').Value = 'BitMartToken'
$ws.Range('C13').Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '0.09036'
$ws.Range('E13').NumberFormat = "@"
$ws.Range('E13').Value = '-0.34%'

# Row 14
$ws.Range('B14').Value = 'BitForexToken'
$ws.Range('C14').Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '0.001571'
$ws.Range('E14').NumberFormat = "@"
$ws.Range('E14').Value = '1.22%'

# Row 15
$ws.Range('B15').Value = 'One'
$ws.Range('C15').Value = 'https://coinranking.com/coin/6Lga5NiXX3rT+one-one'
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '0.0006351'
$ws.Range('E15').NumberFormat = "@"
$ws.Range('E15').Value = '1.52%'

# Row 16
$ws.Range('B16').Value = 'TigerCash'
$ws.Range('C16').Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '0.005778'
$ws.Range('E16').NumberFormat = "@"
$ws.Range('E16').Value = '-3.36%'

# Row 17
$ws.Range('B17').Value = 'LEO'
$ws.Range('C17').Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '3.450'
$ws.Range('E17').NumberFormat = "@"
$ws.Range('E17').Value = '-1.49%'

# Row 18
$ws.Range('B18').Value = 'GateToken'
$ws.Range('C18').Value = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '3.299'
$ws.Range('E18').NumberFormat = "@"
$ws.Range('E18').Value = '-1.25%'

# Row 19
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '2.284'
$ws.Range('E19').NumberFormat = "@"
$ws.Range('E19').Value = '-0.64%'

# Row 21
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '0.1323'
$ws.Range('E21').NumberFormat = "@"
$ws.Range('E21').Value = '2.69%'

# Row 22
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '3.903'
$ws.Range('E22').NumberFormat = "@"
$ws.Range('E22').Value = '-0.99%'

# Row 23
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '0.04400'
$ws.Range('E23').NumberFormat = "@"
$ws.Range('E23').Value = '1.37%'

# Row 24
$ws.Range('E24').NumberFormat = "@"
$ws.Range('E24').Value = '0.08%'

# Row 25
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '0.003841'
$ws.Range('E25').NumberFormat = "@"
$ws.Range('E25').Value = '5.00%'

# Row 26
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '0.0001200'
$ws.Range('E26').NumberFormat = "@"
$ws.Range('E26').Value = '0.26%'

# Row 27
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '0.0001935'
$ws.Range('E27').NumberFormat = "@"
$ws.Range('E27').Value = '14.58%'

# Row 40
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '0.04114'
$ws.Range('E40').NumberFormat = "@"
$ws.Range('E40').Value = '0.87%'

# Row 41
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '0.006771'
$ws.Range('E41').NumberFormat = "@"
$ws.Range('E41').Value = '2.03%'

# Row 42
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '0.1171'
$ws.Range('E42').NumberFormat = "@"
$ws.Range('E42').Value = '-0.35%'

# Row 43
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '0.001941'
$ws.Range('E43').NumberFormat = "@"
$ws.Range('E43').Value = '-17.23%'

# Row 44
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '0.01149'
$ws.Range('E44').NumberFormat = "@"
$ws.Range('E44').Value = '-8.35%'

# Row 45
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '0.00005240'
$ws.Range('E45').NumberFormat = "@"
$ws.Range('E45').Value = '-0.06%'

# Row 46
$ws.Range('B46').Value = 'BOLO'
$ws.Range('C46').Value = 'https://coinranking.com/coin/ogrGe0dEab+bolo-bolo'
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '1.490'
$ws.Range('E46').NumberFormat = "@"
$ws.Range('E46').Value = '-36.58%'

# Row 47
$ws.Range('B47').Value = 'CoinbaseStockToken'
$ws.Range('C47').Value = 'https://coinranking.com/coin/_ZA6fIr53+coinbasestocktoken-coin'
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '0.01999'
$ws.Range('E47').NumberFormat = "@"
$ws.Range('E47').Value = '-11.31%'
